$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow {
    param(
        [int]$Row,
        [string]$Company,
        [string]$Referral,
        [string]$Progress,
        [string]$Style
    )

    $cellA = $ws.Range("A$Row")
    $cellB = $ws.Range("B$Row")
    $cellC = $ws.Range("C$Row")

    $cellA.Value = $Company
    $cellB.Value = $Referral
    $cellC.Value = $Progress

    if ($Style -eq "yellow") {
        # Matches the existing highlighted "interview completed" style (fontId 2 / fillId 2).
        $range = $ws.Range("A$Row" + ":C$Row")
        $range.Font.Name = "微软雅黑"
        $range.HorizontalAlignment = -4108
        $range.Interior.Color = 65535
    }
    elseif ($Style -eq "white") {
        # Row header cell (company) and progress cell keep the plain centered style,
        # while the referral cell gets a new explicit white-fill centered style.
        $cellA.Font.Name = "微软雅黑"
        $cellA.HorizontalAlignment = -4108
        $cellC.Font.Name = "微软雅黑"
        $cellC.HorizontalAlignment = -4108

        $cellB.Font.Name = "微软雅黑"
        $cellB.HorizontalAlignment = -4108
        $cellB.Interior.Color = 16777215
    }
    else {
        $range = $ws.Range("A$Row" + ":C$Row")
        $range.Font.Name = "微软雅黑"
        $range.HorizontalAlignment = -4108
    }
}

# Three repeated blocks: a newly tracked company pair ("招银网络科技" / "腾讯")
# followed by the same three companies that were already being tracked,
# plus one trailing pair-only block (rows 20-21).
for ($block = 0; $block -lt 3; $block++) {
    $base = 5 + ($block * 5)

    Set-DataRow ($base + 0) "招银网络科技" "是" "简历投递" "white"
    Set-DataRow ($base + 1) "腾讯" "否" "简历投递" "normal"
    Set-DataRow ($base + 2) "酷狗" "是" "简历投递" "normal"
    Set-DataRow ($base + 3) "OPPO广东移动通信有限公司" "是" "简历投递" "normal"
    Set-DataRow ($base + 4) "华为" "是" "面试完成" "yellow"
}

Set-DataRow 20 "招银网络科技" "是" "简历投递" "white"
Set-DataRow 21 "腾讯" "否" "简历投递" "normal"

$ws.Columns("C").ColumnWidth = 23

$ws.Range("A17:C21").Select()
